$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to Text format first so numeric-looking strings
    # (e.g. "314.91", "0.07110") are kept as literal text instead of
    # being parsed into a floating point number, then restore the
    # default "Normal" style so no stray number-format style lingers.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "24.646.43"
$ws.Range("E2").Value = "  -0.04%  "
Set-TextValue "D3" "1.692.84"
$ws.Range("E3").Value = "  +0.08%  "
Set-TextValue "D4" "1.005"
$ws.Range("E4").Value = "  +0.29%  "
Set-TextValue "D5" "314.91"
$ws.Range("E5").Value = "  -0.22%  "
Set-TextValue "D6" "0.9993"
$ws.Range("E6").Value = "  -0.27%  "
Set-TextValue "D7" "0.3907"
Set-TextValue "D8" "0.4036"
$ws.Range("E8").Value = "  -0.68%  "
Set-TextValue "D9" "1.492"
$ws.Range("E9").Value = "  -0.13%  "
Set-TextValue "D10" "1.006"
$ws.Range("E10").Value = "  +0.35%  "
Set-TextValue "D11" "52.83"
$ws.Range("E11").Value = "  -0.79%  "
Set-TextValue "D12" "0.08732"
$ws.Range("E12").Value = "  -0.92%  "
Set-TextValue "D13" "7.604"
$ws.Range("E13").Value = "  +5.21%  "
Set-TextValue "D14" "24.73"
$ws.Range("E14").Value = "  +5.04%  "
Set-TextValue "D15" "0.00001353"
$ws.Range("E15").Value = "  +2.87%  "
Set-TextValue "D16" "7.936"
$ws.Range("E16").Value = "  -1.59%  "
Set-TextValue "D17" "1.684.43"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("E18").Value = "  -1.44%  "
Set-TextValue "D19" "0.07110"
$ws.Range("E19").Value = "  +1.36%  "
Set-TextValue "D20" "19.82"
$ws.Range("E20").Value = "  +1.71%  "
Set-TextValue "D21" "7.291"
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  -0.17%  "
Set-TextValue "D24" "24.637.01"
$ws.Range("E24").Value = "  -0.03%  "
Set-TextValue "D25" "3.015"
$ws.Range("E25").Value = "  -7.72%  "
Set-TextValue "D26" "2.350"
$ws.Range("E26").Value = "  -0.48%  "
Set-TextValue "D27" "22.72"
$ws.Range("E27").Value = "  -0.08%  "
Set-TextValue "D28" "161.73"
$ws.Range("E28").Value = "  -0.86%  "
Set-TextValue "D29" "8.372"
$ws.Range("E29").Value = "  +10.21%  "
Set-TextValue "D30" "5.246"
$ws.Range("E30").Value = "  +1.06%  "
Set-TextValue "D31" "136.64"
$ws.Range("E31").Value = "  +0.66%  "
Set-TextValue "D32" "1.869.41"
$ws.Range("E32").Value = "  -0.82%  "
Set-TextValue "D33" "0.08855"
$ws.Range("E33").Value = "  +3.42%  "
Set-TextValue "D34" "7.514"
$ws.Range("E34").Value = "  +4.62%  "
Set-TextValue "D35" "1.038"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  +4.33%  "
Set-TextValue "D37" "0.02926"
$ws.Range("E37").Value = "  +7.58%  "
Set-TextValue "D38" "0.2727"
$ws.Range("E38").Value = "  -0.22%  "
Set-TextValue "D39" "10.72"
$ws.Range("E39").Value = "  -5.50%  "
Set-TextValue "D40" "14.28"
$ws.Range("E40").Value = "  -0.39%  "
Set-TextValue "D41" "0.09116"
$ws.Range("E41").Value = "  -0.83%  "
Set-TextValue "D42" "0.7833"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("E43").Value = "  +0.37%  "
Set-TextValue "D44" "16.77"
$ws.Range("E44").Value = "  +5.31%  "
Set-TextValue "D45" "0.7201"
$ws.Range("E45").Value = "  +0.90%  "
Set-TextValue "D46" "2.578"
Set-TextValue "D47" "4.208"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +0.17%  "
Set-TextValue "D49" "1.327"
$ws.Range("E49").Value = "  +0.63%  "
Set-TextValue "D50" "138.69"
$ws.Range("E50").Value = "  -0.35%  "
Set-TextValue "D51" "91.03"
$ws.Range("E51").Value = "  +1.37%  "
